# "add thank you slide" - turn the big colored "Thank You" headline on the
# last slide into a small plain "Thank you" textbox.
#
# The canonical OOXML diff shows the old shape (id=4 "TextBox 3") being
# replaced by a brand-new shape (id=2 "TextBox 1") at a new position/size
# with plain (un-styled) text - i.e. the original textbox was deleted and a
# fresh one drawn in its place, rather than simply being restyled in place.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(3)

# Remove the old "Thank You" textbox.
$old = $s.Shapes.Item("TextBox 3")
$old.Delete()

# Draw the replacement textbox. On this (now empty) slide the freed-up
# lowest id/name is reused, giving us id=2 / "TextBox 1" to match the diff.
$new = $s.Shapes.AddTextbox(1, 100, 100, 200, 50)
$new.Name = "TextBox 1"

# Set the text before the final size/position so PowerPoint's "resize shape
# to fit text" (spAutoFit) autosizing doesn't get a chance to overwrite the
# exact width/height we set below.
$tr = $new.TextFrame.TextRange
$tr.Text = "Thank you"

$new.TextFrame.WordWrap = $false
$new.TextFrame.AutoSize = 1
$new.Fill.Visible = $false

# Final geometry (EMU shown in comments; Left/Top/Width/Height are in points,
# 1 pt = 12700 EMU):
#   off  x=5340096 y=1664208
#   ext cx=1153393 cy=369332
$new.Left = 420.48
$new.Top = 131.04001
$new.Width = 90.8184
$new.Height = 29.081259842519685
